$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.857.70'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '1.627.14'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'214.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").Value = "'0.501"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("D10").Value = "'19.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.62%  '
$ws.Range("D11").Value = "'0.0788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.68%  '
$ws.Range("D12").Value = '1.852.67'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").Value = "'4.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").Value = '1.634.37'
$ws.Range("E14").Value = '  +0.49%  '
$ws.Range("D15").Value = "'0.542"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.46%  '
$ws.Range("D16").Value = '0.0₃0755'
$ws.Range("E16").Value = '  -0.60%  '
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").Value = '25.852.77'
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = "'192.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.48%  '
$ws.Range("E21").Value = '  -1.85%  '
$ws.Range("D22").Value = "'9.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("D23").Value = "'6.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("E24").Value = '  -1.78%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").Value = "'142.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.24%  '
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("D28").Value = "'6.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").Value = "'15.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").Value = "'0.0497"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.49%  '
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("E33").Value = '  -0.18%  '
$ws.Range("D34").Value = "'1.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.57%  '
$ws.Range("E35").Value = '  +1.12%  '
$ws.Range("D36").Value = "'0.900"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.20%  '
$ws.Range("D37").Value = '1.127.86'
$ws.Range("E37").Value = '  -0.48%  '
$ws.Range("E39").Value = '  -2.24%  '
$ws.Range("E40").Value = '  +0.97%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").Value = "'5.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("E43").Value = '  -0.74%  '
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").Value = '1.763.46'
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("E46").Value = '  -1.03%  '
$ws.Range("D47").Value = "'56.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.38%  '
$ws.Range("D48").Value = "'0.0530"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.42%  '
$ws.Range("E50").Value = '  -0.79%  '
$ws.Range("D51").Value = "'7.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.39%  '
